$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.931.40'
$ws.Range('E2').Value = '  -3.54%  '
$ws.Range('D3').Value = '3.374.94'
$ws.Range('E3').Value = '  -4.50%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.87'
$ws.Range('E5').Value = '  -4.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.00'
$ws.Range('E6').Value = '  -6.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.600'
$ws.Range('E7').Value = '  -2.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '3.368.36'
$ws.Range('E9').Value = '  -4.32%  '
$ws.Range('E10').Value = '  -8.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.596'
$ws.Range('E11').Value = '  -4.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.15'
$ws.Range('E12').Value = '  -7.25%  '
$ws.Range('E13').Value = '  -5.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.77'
$ws.Range('E14').Value = '  -5.95%  '
$ws.Range('D15').Value = '3.910.27'
$ws.Range('E15').Value = '  -4.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '608.03'
$ws.Range('E16').Value = '  -11.58%  '
$ws.Range('D17').Value = '66.854.28'
$ws.Range('E17').Value = '  -3.69%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.381.41'
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.99'
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('E20').Value = '  -2.99%  '
$ws.Range('E21').Value = '  -6.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.917'
$ws.Range('E22').Value = '  -5.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.09'
$ws.Range('E23').Value = '  -4.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.14'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.38'
$ws.Range('E25').Value = '  -9.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.11'
$ws.Range('E26').Value = '  -6.78%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.03'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.77'
$ws.Range('E28').Value = '  -5.82%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.54'
$ws.Range('E29').Value = '  -7.74%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.87'
$ws.Range('E30').Value = '  -8.68%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.87'
$ws.Range('E31').Value = '  -7.94%  '
$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.98'
$ws.Range('E32').Value = '  -9.68%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.36'
$ws.Range('E33').Value = '  -8.03%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.22'
$ws.Range('E34').Value = '  -6.01%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '3.877.69'
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.106'
$ws.Range('E36').Value = '  -5.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '541.06'
$ws.Range('E37').Value = '  +6.63%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '58.48'
$ws.Range('E38').Value = '  -6.48%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.46'
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('E41').Value = '  -11.18%  '
$ws.Range('B42').Value = 'CoreDAO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.45'
$ws.Range('E42').Value = '  +30.64%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.74'
$ws.Range('E43').Value = '  -7.73%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.128'
$ws.Range('E44').Value = '  -5.61%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.352'
$ws.Range('E45').Value = '  -5.51%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.59'
$ws.Range('E46').Value = '  -6.45%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0422'
$ws.Range('E47').Value = '  -7.20%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.21'
$ws.Range('E48').Value = '  -5.48%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.70'
$ws.Range('E49').Value = '  -8.84%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.131'
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.08%  '
